# UPDATE: ans_extractor for parse json, and DEBUG: sql_refine
#
# The "tables" sheet header in column E was "tb_promptLit" (mixed case);
# rename it to "tb_promptlit" (lowercase) to match the refined schema.

$wb = $excel.ActiveWorkbook

$tables = $wb.Worksheets.Item("tables")
$tables.Range("E1").Value = "tb_promptlit"

# Make "tables" the active/selected sheet, with E1 selected.
$tables.Activate()
$tables.Range("E1").Select()
